# Append the 2025-11-19 allocation row (row 79) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A to stay as text so the date-looking string "11/19/2025"
# is not auto-converted into a serial date number, then reset the style
# back to Normal so no stray per-cell formatting is left behind.
$ws.Range("A79").NumberFormat = "@"
$ws.Range("A79").Value = "11/19/2025"
$ws.Range("A79").Style = "Normal"

$ws.Range("B79").Value = 0.2041669959013599
$ws.Range("C79").Value = 0.7958330040986401
